$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the three new localization rows (key / default / ja-JP) in
# alphabetical-key order, shifting everything below them down.

# "copy" goes right before the existing "create" row (269)
$ws.Rows("269:269").Insert()
$ws.Cells.Item(269, 1).Value = "copy"
$ws.Cells.Item(269, 2).Value = "Copy"
$ws.Cells.Item(269, 4).Value = "コピー"

# "cut" goes right before the existing "danger" row (originally 271,
# now at 272 after the previous insert)
$ws.Rows("272:272").Insert()
$ws.Cells.Item(272, 1).Value = "cut"
$ws.Cells.Item(272, 2).Value = "Cut"
$ws.Cells.Item(272, 4).Value = "カット"

# "paste" goes right before the existing "please_select" row (originally
# 378, now at 380 after the two previous inserts)
$ws.Rows("380:380").Insert()
$ws.Cells.Item(380, 1).Value = "paste"
$ws.Cells.Item(380, 2).Value = "Paste"
$ws.Cells.Item(380, 4).Value = "ペースト"

# Match the author's final selection / scroll position.
$ws.Range("M380").Select()
